$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the existing data rows (old row 2 shifts to row 5)
$ws.Rows("2:4").Insert()

# Row 2 (new record, ID 7)
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "'0105781496"
$ws.Cells.Item(2, 3).Value = 45257.41307870371
$ws.Cells.Item(2, 4).Value = "Nyw5LR"

# Row 3 (new record, ID 6)
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = "'0105781496"
$ws.Cells.Item(3, 3).Value = 45257.41237268518
$ws.Cells.Item(3, 4).Value = "'111"

# Row 4 (new record, ID 5)
$ws.Cells.Item(4, 1).Value = 5
$ws.Cells.Item(4, 2).Value = "'0105781496"
$ws.Cells.Item(4, 3).Value = 45257
$ws.Cells.Item(4, 4).Value = "'111"

# Apply date+time format to the whole FECHA column (new rows + pre-existing rows)
$ws.Range("C2:C8").NumberFormat = "yyyy-mm-dd h:mm:ss"

Write-Host "done"
